$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 7
$ws.Range("G4").Value = 8
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 12
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 9
$ws.Range("G9").Value = 5
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 4
$ws.Range("G12").Value = 4
$ws.Range("G13").Value = 5
$ws.Range("G14").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("G16").Value = 9
$ws.Range("G17").Value = 2
$ws.Range("G18").Value = 5
$ws.Range("G19").Value = 6
$ws.Range("G20").Value = 7
$ws.Range("G21").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("G23").Value = 5
$ws.Range("G24").Value = 3
$ws.Range("G25").Value = 5
$ws.Range("G26").Value = 4
$ws.Range("G27").Value = 5
$ws.Range("G28").Value = 2
$ws.Range("G29").Value = 0
